$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the changed values in existing row 54 ---
# (columns whose value did not change are intentionally left untouched)
$row54 = [ordered]@{
    "B"  = 253.1
    "C"  = 291.1
    "D"  = 126.6
    "F"  = 119.1
    "G"  = 62.5
    "H"  = 89.90000000000001
    "I"  = 89.5
    "J"  = 112.9
    "K"  = 36.8
    "L"  = 80.2
    "M"  = 93.40000000000001
    "P"  = 142.7
    "Q"  = 100.3
    "R"  = 123.9
    "S"  = 150.1
    "U"  = 102
    "V"  = 84.09999999999999
    "Y"  = 89.3
    "Z"  = 81.40000000000001
    "AA" = 105.8
    "AC" = 89.3
    "AE" = 91.3
    "AF" = 115.4
    "AG" = 136
    "AH" = 110.4
    "AI" = 103.1
    "AL" = 88.7
    "AM" = 106.2
}

foreach ($col in $row54.Keys) {
    $ws.Range("$col" + "54").Value = $row54[$col]
}

# --- Append new row 55 with the new quarter ---
# "01-04-2021" looks like a date, so a direct .Value assignment would be
# auto-converted into a date serial number by Excel's smart type detection.
# Going through a text formula + paste-values round trip keeps it as a
# literal, unstyled shared string, matching how the other period labels
# in column A are stored.
$dateCell = $ws.Range("A55")
$helper = $ws.Range("ZZ1")
$helper.Formula = '="01-04-2021"'
$helper.Copy()
$dateCell.PasteSpecial(-4163) # xlPasteValues
$helper.Clear()

$row55 = [ordered]@{
    "B"  = 110.2
    "C"  = 114.2
    "D"  = 136.3
    "E"  = 46.5
    "F"  = 125.6
    "G"  = 94.09999999999999
    "H"  = 99.59999999999999
    "I"  = 98.90000000000001
    "J"  = 130.9
    "K"  = 46.3
    "L"  = 80.90000000000001
    "M"  = 105.9
    "N"  = 111.3
    "O"  = 47.8
    "P"  = 161.6
    "Q"  = 105.9
    "R"  = 124.2
    "S"  = 117.4
    "T"  = 117.9
    "U"  = 128.9
    "V"  = 97.59999999999999
    "W"  = 110.3
    "X"  = 42.9
    "Y"  = 98.90000000000001
    "Z"  = 96.7
    "AA" = 105.2
    "AB" = 78
    "AC" = 95
    "AD" = 98
    "AE" = 84.3
    "AF" = 133.8
    "AG" = 156.3
    "AH" = 128.3
    "AI" = 103.5
    "AJ" = 84.3
    "AK" = 59.5
    "AL" = 87.09999999999999
    "AM" = 103.2
}

foreach ($col in $row55.Keys) {
    $ws.Range("$col" + "55").Value = $row55[$col]
}
